$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bring over the font formatting used by the rest of the parameter
# table before filling in values, so the new rows look like their
# neighbours (copy format only, no values) ---

# Rows 19 & 20 mirror row 17's per-column formatting (A-F all share the
# same font treatment there).
$ws.Range("A17:F17").Copy() | Out-Null
$ws.Range("A19:F20").PasteSpecial(-4122) | Out-Null

# Rows 21 & 22 mirror row 18's per-column formatting exactly (names/units
# column + yes/no column on one font, the bmin/bmax/islog/p0 columns on
# the other).
$ws.Range("A18:F18").Copy() | Out-Null
$ws.Range("A21:F21").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:F22").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- New row 19: Silence_LacI_rep ---
$ws.Range("A19").Value = "Silence_LacI_rep"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = "no"
$ws.Range("G19").Value = "silence"

# --- New row 20: pt7_LacI ---
$ws.Range("A20").Value = "pt7_LacI"
$ws.Range("B20").Value = 0.00001
$ws.Range("C20").Value = 100
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 98
$ws.Range("F20").Value = "yes"
$ws.Range("G20").Value = "k_{pt7_PacI}"

# --- Row 21 (previously a blank placeholder row): P3_Lacn_5_cit ---
$ws.Range("A21").Value = "P3_Lacn_5_cit"
$ws.Range("B21").Value = 0.00001
$ws.Range("C21").Value = 100
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 92
$ws.Range("F21").Value = "yes"
$ws.Range("G21").Value = "k_{Cit_Lacn3}"

# --- Row 22 (previously a blank placeholder row): P3_Lacn_5_cit_L ---
$ws.Range("A22").Value = "P3_Lacn_5_cit_L"
$ws.Range("B22").Value = 0.00001
$ws.Range("C22").Value = 0.01
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0.0003
$ws.Range("F22").Value = "yes"
$ws.Range("G22").Value = "kL_cit_pt7"

# --- Update the active selection to B23 (matches the saved cursor state) ---
$ws.Range("B23").Select() | Out-Null

# --- Window geometry tweak recorded in the workbook view ---
$wnd = $excel.ActiveWindow
$wnd.Left = 15680
$wnd.Top = 760
$wnd.Width = 30240
$wnd.Height = 17700
